# Generate Report for Handoff
#
# The handoff report moved from "Handed back: in sync with en-US" to
# "Ready for handoff", and the associated timestamps were refreshed.
# Because the new status text is noticeably shorter than the old one,
# the status columns are narrowed to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-locale status + the "Latest HO Xliff Generate Date" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-28 16:59:53"

# --- zh-cn detail sheet: Status + Latest Handoff Datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-28 16:59:49"

# --- de-de detail sheet: Status + Latest Handoff Datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-28 16:59:53"

# Shrink the status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
